$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entries for 김대환 (row 5) and 신보라 (row 12): clear their
# name/university/address cells but keep the row (and A column's style)
# in place, same as category-based cleanup described in the commit.
$ws.Range("A5:C5").ClearContents()
$ws.Range("A12:C12").ClearContents()

# Update the saved selection/active cell to B3
$ws.Range("B3").Select()
